# DUNEXMainExp_notes.xlsx - "added up to mission 50 to notes"
# Adds mission rows 41-50 (sheet rows 43-52) with deployed/retrieved
# microSWIFT lists and start/end times, matches formatting of existing
# rows, and updates the active sheet/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 43 (Mission 41) ---
$ws.Range("A43").Value = 41
$ws.Range("N43").Value = "49,46,50,44,43,42,37,36,35,34,33,32,31,58,23"
$ws.Range("O43").Value = "49,46,50,44,43,42,37,36,35,34,33,32,31,58,23"
$ws.Range("Q43").Value = "2021-10-17T12:15:00"
$ws.Range("R43").Value = "2021-10-17T13:00:00"
$ws.Rows.Item(43).RowHeight = 34

# --- Row 44 (Mission 42) ---
$ws.Range("A44").Value = 42
$ws.Range("N44").Value = "42,43,44,46,49,50,31,32,34,35,36,37,23,58,56,2,3,4,16,17,10,27,59,19,60,29"
$ws.Range("O44").Value = "42,43,44,46,49,50,31,32,34,35,36,37,23,58,56,2,3,4,16,17,10,27,59,19,60,29"
$ws.Range("Q44").Value = "2021-10-17T13:00:00"
$ws.Range("R44").Value = "2021-10-17T14:00:00"
$ws.Rows.Item(44).RowHeight = 51

# --- Row 45 (Mission 43) ---
$ws.Range("A45").Value = 43
$ws.Range("N45").Value = "46,42,36,2,3,4,35,33,29,17,58,23,10,19,32,27,44,31,49,13,57,14,50,43,56,59,16,34"
$ws.Range("O45").Value = "46,42,36,2,3,4,35,33,29,17,58,23,10,19,32,27,44,31,49,13,57,14,50,43,56,59,16,34"
$ws.Range("Q45").Value = "2021-10-17T15:22:00"
$ws.Range("R45").Value = "2021-10-17T16:30:00"
$ws.Rows.Item(45).RowHeight = 51

# --- Row 46 (Mission 44) ---
$ws.Range("A46").Value = 44
$ws.Range("N46").WrapText = $true
$ws.Range("O46").WrapText = $true
$ws.Range("N46").Value = "4,2,3,4,8,56,9,10,11,12,13,14,16,17,18,19,20,21,23,57"
$ws.Range("O46").Value = "4,2,3,4,8,56,9,10,11,12,13,14,16,17,18,19,20,21,23,57"
$ws.Range("Q46").Value = "2021-10-18T13:09:00"
$ws.Range("R46").Value = "2021-10-18T15:42:00"
$ws.Rows.Item(46).RowHeight = 34

# --- Row 47 (Mission 45) ---
$ws.Range("A47").Value = 45
$ws.Range("N47").Value = "56,2,3,4,8,10,12,13,14,57,16,17,18,19,20,21,23,24,58,42"
$ws.Range("O47").Value = "56,2,3,4,8,10,12,13,14,57,16,17,18,19,20,21,23,24,58,42"
$ws.Range("Q47").Value = "2021-10-18T17:15:00"
$ws.Range("R47").Value = "2021-10-18T18:40:00"
$ws.Rows.Item(47).RowHeight = 34

# --- Row 48 (Mission 46) ---
$ws.Range("A48").Value = 46
$ws.Range("N48").WrapText = $true
$ws.Range("O48").WrapText = $true
$ws.Range("N48").Value = "41,42,43,44,46,49,50,31,32,33,34,35"
$ws.Range("O48").Value = "41,42,43,44,46,49,50,31,32,33,34,35"
$ws.Range("Q48").Value = "2021-10-19T12:03:00"
$ws.Range("R48").Value = "2021-10-19T12:13:00"
$ws.Rows.Item(48).RowHeight = 34

# --- Row 49 (Mission 47) ---
$ws.Range("A49").Value = 47
$ws.Range("N49").Value = "41,42,43,44,46,49,50,31,32,33,34,35"
$ws.Range("O49").Value = "41,42,43,44,46,49,50,31,32,33,34,35"
$ws.Range("Q49").Value = "2021-10-19T12:45:00"
$ws.Range("R49").Value = "2021-10-19T12:50:00"
$ws.Rows.Item(49).RowHeight = 34

# --- Row 50 (Mission 48) ---
$ws.Range("A50").Value = 48
$ws.Range("N50").WrapText = $true
$ws.Range("O50").WrapText = $true
$ws.Range("N50").Value = "41,42,43,44,46,49,50,31,32,33,34,35"
$ws.Range("O50").Value = "41,42,43,44,46,49,50,31,32,33,34,35"
$ws.Range("Q50").Value = "2021-10-19T13:18:00"
$ws.Range("R50").Value = "2021-10-19T13:28:00"
$ws.Rows.Item(50).RowHeight = 34

# --- Row 51 (Mission 49) ---
$ws.Range("A51").Value = 49
$ws.Range("N51").Value = "41,42,43,44,46,49,50,31,32,33,34,35"
$ws.Range("O51").Value = "41,42,43,44,46,49,50,31,32,33,34,35"
$ws.Range("Q51").Value = "2021-10-19T13:40:00"
$ws.Range("R51").Value = "2021-10-19T13:50:00"
$ws.Rows.Item(51).RowHeight = 34

# --- Row 52 (Mission 50) ---
$ws.Range("A52").Value = 50
$ws.Range("N52").WrapText = $true
$ws.Range("O52").WrapText = $true
$ws.Range("N52").Value = "31,33,34,35,36,37,21,23,24,58,27,13,14,2,4,57,16,17,9,56,20,29,8"
$ws.Range("O52").Value = "31,33,34,35,36,37,21,23,24,58,27,13,14,2,4,57,16,17,9,56,20,29,8"
$ws.Range("Q52").Value = "2021-10-19T15:15:00"
$ws.Range("R52").Value = "2021-10-19T15:36:00"
$ws.Rows.Item(52).RowHeight = 51

# --- Carry the "microSWIFTs Deployed/Retrieved" wrap format further down
# (rows 54-74, every other row) to match the formatting paste that was
# applied across the rest of the blank mission rows below. ---
$ws.Range("N43:O43").Copy()
$ws.Range("N54:O74").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the active sheet / view / selection to match the saved state ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("C39").Select()

$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("R53").Select()
